{"js": "// The commit adds a short sentence to the end of the paragraph that\n// discusses earlier censuses' \"third gender\" option (Nepal/India/Pakistan).\n// Find that paragraph by its distinctive trailing text and append the new\n// run \" Louise edited this line. \" after it \u2014 this mirrors the OOXML diff,\n// which adds a brand-new <w:r><w:t> run at the end of the paragraph,\n// immediately after \"...with less detailed reporting of results.\"\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"with less detailed reporting of results.\";\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (text && text.indexOf(marker) !== -1) {\n    target = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not locate target paragraph ending in: \" + marker);\n}\n\n// Insert a new run of text at the end of the paragraph (after the existing\n// text, before the paragraph mark) \u2014 matches the new <w:r> appended in the\n// diff.\ntarget.insertText(\" Louise edited this line. \", Word.InsertLocation.end);\n\nawait context.sync();\n", "ps1": "# The commit appends a short new sentence to the end of the paragraph that\n# discusses earlier censuses' \"third gender\" option (Nepal/India/Pakistan).\n# In the OOXML diff this shows up as a brand-new <w:r><w:t> run added right\n# after \"...with less detailed reporting of results.\" and before the\n# closing </w:p>.\n\n$d = $word.ActiveDocument\n\n$target = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    if ($p.Range.Text -like \"*with less detailed reporting of results.*\") {\n        $target = $p\n        break\n    }\n}\n\nif ($null -eq $target) {\n    throw \"Could not locate target paragraph ending in 'with less detailed reporting of results.'\"\n}\n\n# Append the new run of text at the end of the paragraph's range (i.e.\n# right before the paragraph mark) -- matches the new run added in the diff.\n$target.Range.InsertAfter(\" Louise edited this line. \")\n"}
